# Update countries & provincias Spain
# Applies the refreshed COVID-19 "paises" data snapshot (19 May 2020, 05:05)
# on top of the previous snapshot (04:35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 05:05"

# --- Update case counters for countries whose figures changed but whose ---
# --- rank (sort order) in the table stayed the same                     ---

# Row 60: Kazajistan
$ws.Cells.Item(60, 2).Value = 6751   # Casos totales
$ws.Cells.Item(60, 3).Value = 311    # Nuevos casos
$ws.Cells.Item(60, 5).Value = 3247   # Recuperados

# Row 84: Guatemala
$ws.Cells.Item(84, 2).Value = 2001   # Casos totales
$ws.Cells.Item(84, 3).Value = 89     # Nuevos casos
$ws.Cells.Item(84, 4).Value = 139    # Casos activos
$ws.Cells.Item(84, 5).Value = 1824   # Recuperados
$ws.Cells.Item(84, 7).Value = 3      # Muertes hoy
$ws.Cells.Item(84, 8).Value = 38     # Muertes

# Row 93: Nueva Zelanda
$ws.Cells.Item(93, 2).Value = 1503   # Casos totales
$ws.Cells.Item(93, 3).Value = 4      # Nuevos casos
$ws.Cells.Item(93, 4).Value = 1442   # Casos activos
$ws.Cells.Item(93, 5).Value = 40     # Recuperados

# --- Re-order three mutually-tied countries (Nueva Caledonia / Belice / ---
# --- Santa Lucia, all on 18 total cases) into their new relative order: ---
# --- Belice, Santa Lucia, Nueva Caledonia - each row keeps the figures  ---
# --- belonging to the country it now displays.                         ---

$ws.Cells.Item(195, 1).Value = "Belice"
$ws.Cells.Item(195, 4).Value = 16
$ws.Cells.Item(195, 8).Value = 2

$ws.Cells.Item(196, 1).Value = "Santa Lucia"
$ws.Cells.Item(196, 4).Value = 18
$ws.Cells.Item(196, 8).Value = 0

$ws.Cells.Item(197, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(197, 4).Value = 18
$ws.Cells.Item(197, 8).Value = 0

# --- Re-order three mutually-tied countries (Seychelles / Groenlandia / ---
# --- Montserrat, all on 11 total cases) into their new relative order:  ---
# --- Montserrat, Seychelles, Groenlandia - each row keeps the figures   ---
# --- belonging to the country it now displays.                         ---

$ws.Cells.Item(209, 1).Value = "Montserrat"
$ws.Cells.Item(209, 4).Value = 10
$ws.Cells.Item(209, 8).Value = 1

$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Groenlandia"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0
